# Draft mapping update for ror-healthcareservice-contact-telecom
# - bump the IG "Date" metadata value
# - lowercase the top-level "Telecommunication" mapping label
# - clear the business-mapping ("Mapping: ...") annotations that used to sit on
#   the ror-telecom-communication-channel / ror-telecom-usage /
#   ror-telecom-confidentiality-level / telecomAddress extension rows
# - move the "adresseTelecom" mapping down onto the telecomAddress.value[x] row

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

$elements = $wb.Worksheets.Item("Elements")

# Row 2 = the root "Extension" slice -> AL column holds the business mapping
$elements.Range("AL2").Value = "telecommunication"

# Row 5 = Extension.extension:ror-telecom-communication-channel
$elements.Range("AL5").Value = ""

# Row 6 = Extension.extension:ror-telecom-usage
$elements.Range("AL6").Value = ""

# Row 7 = Extension.extension:ror-telecom-confidentiality-level
$elements.Range("AL7").Value = ""

# Row 8 = Extension.extension:telecomAddress
$elements.Range("AL8").Value = ""

# Row 12 = Extension.extension:telecomAddress.value[x] now carries the mapping
$elements.Range("AL12").Value = "adresseTelecom"
